# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" right after the existing "2021-Q4" sheet,
#   holding the fund-holding detail for that quarter (same layout as "2021-Q4").
# - Insert a new leading row into the "总计" (totals) sheet for "2022-Q1",
#   pushing the existing "2021-Q4" total row down.

$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item(1)     # "2021-Q4"

# ------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet, right after "2021-Q4",
#    cloning the layout/formatting of the "2021-Q4" sheet.
# ------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

# NOTE: fetch the "总计" sheet only now (by its new position, 3), since
# inserting a worksheet shifts what Worksheets.Item(2) resolves to.
$wsTotal = $wb.Worksheets.Item(3)  # "总计"

# (copy header row and data row separately so we don't manufacture a
# spurious/empty A1 cell that never existed on the source sheet)
$wsQ4.Range("B1:H1").Copy($wsQ1.Range("B1"))
$wsQ4.Range("A2:H2").Copy($wsQ1.Range("A2"))

# Fill in the 2022-Q1 fund-holding data (row 2). Columns that hold
# numeric-looking text (fund code / scale / position figures) must be
# force-typed as text so values like the leading-zero fund code are not
# coerced into numbers; ClearFormats() afterwards drops the temporary
# "Text" number-format style so the cell ends up unstyled, like the source.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $wsQ1.Range("B2") "005433"
Set-TextValue $wsQ1.Range("C2") "申万菱信医药先锋股票"
Set-TextValue $wsQ1.Range("D2") "2.20"
Set-TextValue $wsQ1.Range("E2") "90.81"
Set-TextValue $wsQ1.Range("F2") "4.80"
Set-TextValue $wsQ1.Range("G2") "0.1056"
$wsQ1.Range("H2").Value = 6

# ------------------------------------------------------------------
# 2) Add a new leading row to the "总计" sheet for "2022-Q1".
# ------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

# Re-apply the index-column style (bold/bordered) to the new row's A cell,
# taken from the row that was just pushed down (still carries that style).
$wsTotal.Range("A3").Copy($wsTotal.Range("A2"))

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.11
